$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.238.47'
$ws.Range("E2").Value = '  +2.31%  '

$ws.Range("D3").Value = '3.394.50'
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.61'
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.06'
$ws.Range("E6").Value = '  +3.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +1.35%  '

$ws.Range("E9").Value = '  +8.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.593'
$ws.Range("E10").Value = '  +2.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.64'
$ws.Range("E11").Value = '  +4.09%  '

$ws.Range("E12").Value = '  +4.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '681.93'
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("E14").Value = '  +3.20%  '

$ws.Range("D15").Value = '3.942.37'
$ws.Range("E15").Value = '  +1.95%  '

$ws.Range("D16").Value = '69.407.51'
$ws.Range("E16").Value = '  +2.56%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.120'
$ws.Range("E17").Value = '  +1.56%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.390.93'
$ws.Range("E18").Value = '  +2.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.76'
$ws.Range("E19").Value = '  +1.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.32'
$ws.Range("E20").Value = '  +2.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.908'
$ws.Range("E21").Value = '  +1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.39'
$ws.Range("E22").Value = '  -2.56%  '

$ws.Range("E23").Value = '  +2.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.98'
$ws.Range("E24").Value = '  +1.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.93'
$ws.Range("E25").Value = '  +0.91%  '

$ws.Range("E26").Value = '  +2.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.73'
$ws.Range("E27").Value = '  +3.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.99'
$ws.Range("E28").Value = '  +3.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.81'
$ws.Range("E29").Value = '  +3.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.98'
$ws.Range("E30").Value = '  -0.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.15'
$ws.Range("E31").Value = '  +1.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '556.99'
$ws.Range("E32").Value = '  -2.19%  '

$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.61'
$ws.Range("E33").Value = '  +10.24%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.54'
$ws.Range("E35").Value = '  +2.28%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").Value = '3.674.00'
$ws.Range("E37").Value = '  -0.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.72'
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("D40").Value = '0.0₃0721'
$ws.Range("E40").Value = '  +7.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.27'
$ws.Range("E41").Value = '  +3.76%  '

$ws.Range("E42").Value = '  +3.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.339'
$ws.Range("E43").Value = '  +1.49%  '

$ws.Range("E44").Value = '  +4.78%  '

$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.69'
$ws.Range("E46").Value = '  +2.04%  '

$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("E48").Value = '  +4.96%  '

$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.03'
$ws.Range("E50").Value = '  +1.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.65'
$ws.Range("E51").Value = '  +2.53%  '
